{"js": "// Replace the worksheet date and every \"a\u00f7b=c, d\" answer cell with its\n// updated value. Each \"before\" string below is unique in the document, so\n// a plain body.search() + insertText(\"Replace\") round-trip is unambiguous.\nconst replacements = [\n  [\"2025-07-11 Friday\", \"2025-07-12 Saturday\"],\n  [\"33\u00f75=6, 3\", \"75\u00f76=12, 3\"],\n  [\"92\u00f75=18, 2\", \"41\u00f75=8, 1\"],\n  [\"72\u00f74=18, 0\", \"36\u00f76=6, 0\"],\n  [\"24\u00f78=3, 0\", \"75\u00f77=10, 5\"],\n  [\"66\u00f72=33, 0\", \"48\u00f78=6, 0\"],\n  [\"57\u00f74=14, 1\", \"82\u00f74=20, 2\"],\n  [\"75\u00f74=18, 3\", \"61\u00f72=30, 1\"],\n  [\"89\u00f74=22, 1\", \"37\u00f72=18, 1\"],\n  [\"13\u00f78=1, 5\", \"28\u00f77=4, 0\"],\n  [\"68\u00f74=17, 0\", \"64\u00f76=10, 4\"],\n  [\"25\u00f76=4, 1\", \"13\u00f73=4, 1\"],\n  [\"71\u00f78=8, 7\", \"15\u00f74=3, 3\"],\n  [\"16\u00f73=5, 1\", \"44\u00f72=22, 0\"],\n  [\"73\u00f72=36, 1\", \"55\u00f76=9, 1\"],\n  [\"58\u00f75=11, 3\", \"32\u00f76=5, 2\"],\n  [\"64\u00f77=9, 1\", \"68\u00f73=22, 2\"],\n  [\"66\u00f78=8, 2\", \"94\u00f77=13, 3\"],\n  [\"51\u00f77=7, 2\", \"70\u00f72=35, 0\"],\n  [\"75\u00f73=25, 0\", \"80\u00f74=20, 0\"],\n  [\"14\u00f76=2, 2\", \"20\u00f73=6, 2\"],\n  [\"42\u00f74=10, 2\", \"31\u00f77=4, 3\"],\n  [\"54\u00f77=7, 5\", \"57\u00f73=19, 0\"],\n  [\"53\u00f76=8, 5\", \"81\u00f73=27, 0\"],\n  [\"71\u00f77=10, 1\", \"47\u00f75=9, 2\"],\n  [\"64\u00f73=21, 1\", \"84\u00f73=28, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the worksheet date and every \"a\u00f7b=c, d\" answer cell with its\n# updated value. Each \"before\" string is unique in the document, so a\n# Find/Replace pass (one per pair, scoped to the whole story) is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-07-11 Friday\", \"2025-07-12 Saturday\"),\n    @(\"33\u00f75=6, 3\", \"75\u00f76=12, 3\"),\n    @(\"92\u00f75=18, 2\", \"41\u00f75=8, 1\"),\n    @(\"72\u00f74=18, 0\", \"36\u00f76=6, 0\"),\n    @(\"24\u00f78=3, 0\", \"75\u00f77=10, 5\"),\n    @(\"66\u00f72=33, 0\", \"48\u00f78=6, 0\"),\n    @(\"57\u00f74=14, 1\", \"82\u00f74=20, 2\"),\n    @(\"75\u00f74=18, 3\", \"61\u00f72=30, 1\"),\n    @(\"89\u00f74=22, 1\", \"37\u00f72=18, 1\"),\n    @(\"13\u00f78=1, 5\", \"28\u00f77=4, 0\"),\n    @(\"68\u00f74=17, 0\", \"64\u00f76=10, 4\"),\n    @(\"25\u00f76=4, 1\", \"13\u00f73=4, 1\"),\n    @(\"71\u00f78=8, 7\", \"15\u00f74=3, 3\"),\n    @(\"16\u00f73=5, 1\", \"44\u00f72=22, 0\"),\n    @(\"73\u00f72=36, 1\", \"55\u00f76=9, 1\"),\n    @(\"58\u00f75=11, 3\", \"32\u00f76=5, 2\"),\n    @(\"64\u00f77=9, 1\", \"68\u00f73=22, 2\"),\n    @(\"66\u00f78=8, 2\", \"94\u00f77=13, 3\"),\n    @(\"51\u00f77=7, 2\", \"70\u00f72=35, 0\"),\n    @(\"75\u00f73=25, 0\", \"80\u00f74=20, 0\"),\n    @(\"14\u00f76=2, 2\", \"20\u00f73=6, 2\"),\n    @(\"42\u00f74=10, 2\", \"31\u00f77=4, 3\"),\n    @(\"54\u00f77=7, 5\", \"57\u00f73=19, 0\"),\n    @(\"53\u00f76=8, 5\", \"81\u00f73=27, 0\"),\n    @(\"71\u00f77=10, 1\", \"47\u00f75=9, 2\"),\n    @(\"64\u00f73=21, 1\", \"84\u00f73=28, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $found = $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
